$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 933028.3  # H9: 1066314.2 -> 933028.3
$ws.Cells.Item(9, 9).Value = 1357042.9  # I9: 1492746.8 -> 1357042.9
$ws.Cells.Item(9, 10).Value = 196.4  # J9: 233 -> 196.4
$ws.Cells.Item(9, 11).Value = 1357042.9  # K9: 1492746.8 -> 1357042.9
$ws.Cells.Item(9, 12).Value = 196.4  # L9: 233 -> 196.4
$ws.Cells.Item(9, 13).Value = -1356873.9  # M9: -1492577.8 -> -1356873.9
$ws.Cells.Item(9, 14).Value = -534.4  # N9: -571 -> -534.4

$ws.Cells.Item(29, 8).Value = 1685.6428  # H29: 1382.0714 -> 1685.6428
$ws.Cells.Item(29, 10).Value = 1807.6154  # J29: 1480.6923 -> 1807.6154
$ws.Cells.Item(29, 12).Value = 5422.8462  # L29: 4442.0769 -> 5422.8462
$ws.Cells.Item(29, 14).Value = -5984.8462  # N29: -5004.0769 -> -5984.8462

$ws.Cells.Item(38, 8).Value = 5499.75  # H38: 4584.7144 -> 5499.75
$ws.Cells.Item(38, 9).Value = 5999  # I38: 3047 -> 5999
$ws.Cells.Item(38, 10).Value = 5333.3335  # J38: 5199.8 -> 5333.3335
$ws.Cells.Item(38, 11).Value = 17997  # K38: 9141 -> 17997
$ws.Cells.Item(38, 12).Value = 16000.0005  # L38: 15599.4 -> 16000.0005
$ws.Cells.Item(38, 13).Value = -17625  # M38: -8769 -> -17625
$ws.Cells.Item(38, 14).Value = -16744.0005  # N38: -16343.4 -> -16744.0005

$ws.Cells.Item(53, 8).Value = 1027.1875  # H53: 1095.5333 -> 1027.1875
$ws.Cells.Item(53, 10).Value = 385.75  # J53: 513.6667 -> 385.75
$ws.Cells.Item(53, 12).Value = 385.75  # L53: 513.6667 -> 385.75
$ws.Cells.Item(53, 14).Value = -1659.75  # N53: -1787.6667 -> -1659.75

$ws.Cells.Item(58, 8).Value = 67  # H58: 64.14286 -> 67
$ws.Cells.Item(58, 9).Value = 67  # I58: 64.14286 -> 67
$ws.Cells.Item(58, 11).Value = 201  # K58: 192.42858 -> 201
$ws.Cells.Item(58, 13).Value = -51  # M58: -42.42858000000001 -> -51

$ws.Cells.Item(112, 8).Value = 1656.421  # H112: 1651.4762 -> 1656.421
$ws.Cells.Item(112, 9).Value = 1299.6666  # I112: 1349.5 -> 1299.6666
$ws.Cells.Item(112, 10).Value = 1723.3125  # J112: 1683.2632 -> 1723.3125
$ws.Cells.Item(112, 11).Value = 3898.9998  # K112: 4048.5 -> 3898.9998
$ws.Cells.Item(112, 12).Value = 5169.9375  # L112: 5049.7896 -> 5169.9375
$ws.Cells.Item(112, 13).Value = -2790.9998  # M112: -2940.5 -> -2790.9998
$ws.Cells.Item(112, 14).Value = -7385.9375  # N112: -7265.7896 -> -7385.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2479.2  # H2: 2549.25 -> 2479.2
$ws.Cells.Item(2, 9).Value = 2479.2  # I2: 2549.25 -> 2479.2
$ws.Cells.Item(2, 11).Value = 2479.2  # K2: 2549.25 -> 2479.2
$ws.Cells.Item(2, 13).Value = -2366.2  # M2: -2436.25 -> -2366.2

$ws.Cells.Item(32, 8).Value = 1436360.9  # H32: 1461124.2 -> 1436360.9
$ws.Cells.Item(32, 9).Value = 662933.2  # I32: 674984.9 -> 662933.2
$ws.Cells.Item(32, 11).Value = 662933.2  # K32: 674984.9 -> 662933.2
$ws.Cells.Item(32, 13).Value = -662646.2  # M32: -674697.9 -> -662646.2

$ws.Cells.Item(37, 8).Value = 49019  # H37: 58037 -> 49019
$ws.Cells.Item(37, 9).Value = 0  # I37: 58037 -> 0
$ws.Cells.Item(37, 10).Value = 49019  # J37: 0 -> 49019
$ws.Cells.Item(37, 11).Value = 0  # K37: 58037 -> 0
$ws.Cells.Item(37, 12).Value = ""  # clear L37 (was 0)
$ws.Cells.Item(37, 13).Value = 49019  # M37: -57764 -> 49019
$ws.Cells.Item(37, 14).Value = -49565  # N37: None -> -49565

$ws.Cells.Item(80, 8).Value = 69403.336  # H80: 59552.25 -> 69403.336
$ws.Cells.Item(80, 10).Value = 94055  # J80: 72703 -> 94055
$ws.Cells.Item(80, 12).Value = 94055  # L80: 72703 -> 94055
$ws.Cells.Item(80, 14).Value = -96051  # N80: -74699 -> -96051

$ws.Cells.Item(83, 8).Value = 69403.336  # H83: 59552.25 -> 69403.336
$ws.Cells.Item(83, 10).Value = 94055  # J83: 72703 -> 94055
$ws.Cells.Item(83, 12).Value = 282165  # L83: 218109 -> 282165
$ws.Cells.Item(83, 14).Value = -292149  # N83: -228093 -> -292149

$ws.Cells.Item(97, 8).Value = 826.6667  # H97: 784.3125 -> 826.6667
$ws.Cells.Item(97, 9).Value = 723  # I97: 678.8461 -> 723
$ws.Cells.Item(97, 11).Value = 723  # K97: 678.8461 -> 723
$ws.Cells.Item(97, 13).Value = -227  # M97: -182.8461 -> -227

$ws.Cells.Item(102, 8).Value = 1938.0454  # H102: 1941.6364 -> 1938.0454
$ws.Cells.Item(102, 9).Value = 1302.2354  # I102: 1306.8823 -> 1302.2354
$ws.Cells.Item(102, 11).Value = 1302.2354  # K102: 1306.8823 -> 1302.2354
$ws.Cells.Item(102, 13).Value = 319.7646  # M102: 315.1177 -> 319.7646

$ws.Cells.Item(116, 8).Value = 2479.2  # H116: 2549.25 -> 2479.2
$ws.Cells.Item(116, 9).Value = 2479.2  # I116: 2549.25 -> 2479.2
$ws.Cells.Item(116, 11).Value = 2479.2  # K116: 2549.25 -> 2479.2
$ws.Cells.Item(116, 13).Value = -185.1999999999998  # M116: -255.25 -> -185.1999999999998

$ws.Cells.Item(122, 8).Value = 2831.5  # H122: 2598.5334 -> 2831.5
$ws.Cells.Item(122, 9).Value = 2777.9  # I122: 2521.4614 -> 2777.9
$ws.Cells.Item(122, 11).Value = 8333.700000000001  # K122: 7564.3842 -> 8333.700000000001
$ws.Cells.Item(122, 13).Value = -5883.700000000001  # M122: -5114.3842 -> -5883.700000000001

$ws.Cells.Item(125, 8).Value = 199998  # H125: 0 -> 199998
$ws.Cells.Item(125, 10).Value = 199998  # J125: 0 -> 199998
$ws.Cells.Item(125, 12).Value = 199998  # L125: 0 -> 199998
$ws.Cells.Item(125, 14).Value = -209838  # N125: None -> -209838

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2479.2  # H3: 2549.25 -> 2479.2
$ws.Cells.Item(3, 9).Value = 2479.2  # I3: 2549.25 -> 2479.2
$ws.Cells.Item(3, 11).Value = 2479.2  # K3: 2549.25 -> 2479.2
$ws.Cells.Item(3, 13).Value = -2365.2  # M3: -2435.25 -> -2365.2

$ws.Cells.Item(64, 8).Value = 1588.3334  # H64: 1430.3 -> 1588.3334
$ws.Cells.Item(64, 9).Value = 1400  # I64: 1350.75 -> 1400
$ws.Cells.Item(64, 10).Value = 2247.5  # J64: 1748.5 -> 2247.5
$ws.Cells.Item(64, 11).Value = 1400  # K64: 1350.75 -> 1400
$ws.Cells.Item(64, 12).Value = 2247.5  # L64: 1748.5 -> 2247.5
$ws.Cells.Item(64, 13).Value = -1175  # M64: -1125.75 -> -1175
$ws.Cells.Item(64, 14).Value = -2697.5  # N64: -2198.5 -> -2697.5

$ws.Cells.Item(67, 8).Value = 1588.3334  # H67: 1430.3 -> 1588.3334
$ws.Cells.Item(67, 9).Value = 1400  # I67: 1350.75 -> 1400
$ws.Cells.Item(67, 10).Value = 2247.5  # J67: 1748.5 -> 2247.5
$ws.Cells.Item(67, 11).Value = 1400  # K67: 1350.75 -> 1400
$ws.Cells.Item(67, 12).Value = 2247.5  # L67: 1748.5 -> 2247.5
$ws.Cells.Item(67, 13).Value = -620  # M67: -570.75 -> -620
$ws.Cells.Item(67, 14).Value = -3807.5  # N67: -3308.5 -> -3807.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 58083.223  # H4: 50236.79 -> 58083.223
$ws.Cells.Item(4, 9).Value = 66666  # I4: 50499.5 -> 66666
$ws.Cells.Item(4, 10).Value = 56366.668  # J4: 50205.883 -> 56366.668
$ws.Cells.Item(4, 11).Value = 66666  # K4: 50499.5 -> 66666
$ws.Cells.Item(4, 12).Value = 56366.668  # L4: 50205.883 -> 56366.668
$ws.Cells.Item(4, 13).Value = -66554  # M4: -50387.5 -> -66554
$ws.Cells.Item(4, 14).Value = -56590.668  # N4: -50429.883 -> -56590.668

$ws.Cells.Item(31, 8).Value = 13891645  # H31: 17859972 -> 13891645
$ws.Cells.Item(31, 9).Value = 998.3333  # I31: 997.5 -> 998.3333
$ws.Cells.Item(31, 10).Value = 20836968  # J31: 25003562 -> 20836968
$ws.Cells.Item(31, 11).Value = 998.3333  # K31: 997.5 -> 998.3333
$ws.Cells.Item(31, 12).Value = 20836968  # L31: 25003562 -> 20836968
$ws.Cells.Item(31, 13).Value = -703.3333  # M31: -702.5 -> -703.3333
$ws.Cells.Item(31, 14).Value = -20837558  # N31: -25004152 -> -20837558

$ws.Cells.Item(34, 8).Value = 13891645  # H34: 17859972 -> 13891645
$ws.Cells.Item(34, 9).Value = 998.3333  # I34: 997.5 -> 998.3333
$ws.Cells.Item(34, 10).Value = 20836968  # J34: 25003562 -> 20836968
$ws.Cells.Item(34, 11).Value = 998.3333  # K34: 997.5 -> 998.3333
$ws.Cells.Item(34, 12).Value = 20836968  # L34: 25003562 -> 20836968
$ws.Cells.Item(34, 13).Value = -796.3333  # M34: -795.5 -> -796.3333
$ws.Cells.Item(34, 14).Value = -20837372  # N34: -25003966 -> -20837372

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 354.0435  # H2: 303.81482 -> 354.0435
$ws.Cells.Item(2, 10).Value = 243.5  # J2: 167.33333 -> 243.5
$ws.Cells.Item(2, 12).Value = 1461  # L2: 1003.99998 -> 1461
$ws.Cells.Item(2, 14).Value = -1687  # N2: -1229.99998 -> -1687

$ws.Cells.Item(33, 8).Value = 257.82352  # H33: 258.375 -> 257.82352
$ws.Cells.Item(33, 10).Value = 259.8  # J33: 262.5 -> 259.8
$ws.Cells.Item(33, 12).Value = 1558.8  # L33: 1575 -> 1558.8
$ws.Cells.Item(33, 14).Value = -2124.8  # N33: -2141 -> -2124.8

$ws.Cells.Item(35, 8).Value = 1529  # H35: 666.5 -> 1529
$ws.Cells.Item(35, 9).Value = 749.75  # I35: 666.5 -> 749.75
$ws.Cells.Item(35, 10).Value = 2152.4  # J35: 0 -> 2152.4
$ws.Cells.Item(35, 11).Value = 2249.25  # K35: 1999.5 -> 2249.25
$ws.Cells.Item(35, 12).Value = 6457.200000000001  # L35: 0 -> 6457.200000000001
$ws.Cells.Item(35, 13).Value = -1961.25  # M35: -1711.5 -> -1961.25
$ws.Cells.Item(35, 14).Value = -7033.200000000001  # N35: None -> -7033.200000000001

$ws.Cells.Item(36, 8).Value = 2700  # H36: 0 -> 2700
$ws.Cells.Item(36, 10).Value = 2700  # J36: 0 -> 2700
$ws.Cells.Item(36, 12).Value = 8100  # L36: 0 -> 8100
$ws.Cells.Item(36, 14).Value = -8438  # N36: None -> -8438

$ws.Cells.Item(39, 8).Value = 9805.833000000001  # H39: 12609.5 -> 9805.833000000001
$ws.Cells.Item(39, 9).Value = 2500  # I39: 0 -> 2500
$ws.Cells.Item(39, 10).Value = 11267  # J39: 12609.5 -> 11267
$ws.Cells.Item(39, 11).Value = 7500  # K39: 0 -> 7500
$ws.Cells.Item(39, 12).Value = 33801  # L39: 37828.5 -> 33801
$ws.Cells.Item(39, 13).Value = -7206  # M39: None -> -7206
$ws.Cells.Item(39, 14).Value = -34389  # N39: -38416.5 -> -34389

$ws.Cells.Item(46, 8).Value = 519.8570999999999  # H46: 560.8570999999999 -> 519.8570999999999
$ws.Cells.Item(46, 9).Value = 317.66666  # I46: 413.33334 -> 317.66666
$ws.Cells.Item(46, 11).Value = 952.9999799999999  # K46: 1240.00002 -> 952.9999799999999
$ws.Cells.Item(46, 13).Value = -861.9999799999999  # M46: -1149.00002 -> -861.9999799999999

$ws.Cells.Item(56, 8).Value = 7665.4  # H56: 8398.5 -> 7665.4
$ws.Cells.Item(56, 9).Value = 7665.4  # I56: 8398.5 -> 7665.4
$ws.Cells.Item(56, 11).Value = 7665.4  # K56: 8398.5 -> 7665.4
$ws.Cells.Item(56, 13).Value = -7135.4  # M56: -7868.5 -> -7135.4

$ws.Cells.Item(128, 8).Value = 175854.72  # H128: 176165.5 -> 175854.72
$ws.Cells.Item(128, 9).Value = 175854.72  # I128: 176165.5 -> 175854.72
$ws.Cells.Item(128, 11).Value = 527564.16  # K128: 528496.5 -> 527564.16
$ws.Cells.Item(128, 13).Value = -522584.16  # M128: -523516.5 -> -522584.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 34999  # H62: 49999 -> 34999
$ws.Cells.Item(62, 10).Value = 34999  # J62: 49999 -> 34999
$ws.Cells.Item(62, 12).Value = 34999  # L62: 49999 -> 34999
$ws.Cells.Item(62, 14).Value = -36371  # N62: -51371 -> -36371

$ws.Cells.Item(65, 8).Value = 34999  # H65: 49999 -> 34999
$ws.Cells.Item(65, 10).Value = 34999  # J65: 49999 -> 34999
$ws.Cells.Item(65, 12).Value = 104997  # L65: 149997 -> 104997
$ws.Cells.Item(65, 14).Value = -111861  # N65: -156861 -> -111861

$ws.Cells.Item(93, 8).Value = 0  # H93: 35000 -> 0
$ws.Cells.Item(93, 10).Value = 0  # J93: 35000 -> 0
$ws.Cells.Item(93, 12).Value = ""  # clear L93 (was 35000)
$ws.Cells.Item(93, 14).Value = 0  # N93: -38744 -> 0

$ws.Cells.Item(102, 8).Value = 4065.8108  # H102: 4257.486 -> 4065.8108
$ws.Cells.Item(102, 9).Value = 1093.4  # I102: 1188.875 -> 1093.4
$ws.Cells.Item(102, 11).Value = 1093.4  # K102: 1188.875 -> 1093.4
$ws.Cells.Item(102, 13).Value = 528.5999999999999  # M102: 433.125 -> 528.5999999999999

$ws.Cells.Item(122, 8).Value = 3080409.5  # H122: 3499791.5 -> 3080409.5
$ws.Cells.Item(122, 9).Value = 4810489  # I122: 4810501 -> 4810489
$ws.Cells.Item(122, 10).Value = 4713.222  # J122: 4565.8335 -> 4713.222
$ws.Cells.Item(122, 11).Value = 14431467  # K122: 14431503 -> 14431467
$ws.Cells.Item(122, 12).Value = 14139.666  # L122: 13697.5005 -> 14139.666
$ws.Cells.Item(122, 13).Value = -14429017  # M122: -14429053 -> -14429017
$ws.Cells.Item(122, 14).Value = -19039.666  # N122: -18597.5005 -> -19039.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5385.5713  # H7: 5983.1665 -> 5385.5713
$ws.Cells.Item(7, 9).Value = 5850  # I7: 9900 -> 5850
$ws.Cells.Item(7, 11).Value = 5850  # K7: 9900 -> 5850
$ws.Cells.Item(7, 13).Value = -5738  # M7: -9788 -> -5738

$ws.Cells.Item(40, 8).Value = 37600.57  # H40: 37486.855 -> 37600.57
$ws.Cells.Item(40, 9).Value = 38867.332  # I40: 38734.668 -> 38867.332
$ws.Cells.Item(40, 11).Value = 38867.332  # K40: 38734.668 -> 38867.332
$ws.Cells.Item(40, 13).Value = -38731.332  # M40: -38598.668 -> -38731.332

$ws.Cells.Item(87, 8).Value = 40000  # H87: 0 -> 40000
$ws.Cells.Item(87, 10).Value = 40000  # J87: 0 -> 40000
$ws.Cells.Item(87, 12).Value = 40000  # L87: 0 -> 40000
$ws.Cells.Item(87, 14).Value = -42246  # N87: None -> -42246

$ws.Cells.Item(90, 8).Value = 40000  # H90: 0 -> 40000
$ws.Cells.Item(90, 10).Value = 40000  # J90: 0 -> 40000
$ws.Cells.Item(90, 12).Value = 120000  # L90: 0 -> 120000
$ws.Cells.Item(90, 14).Value = -131232  # N90: None -> -131232

$ws.Cells.Item(122, 8).Value = 3190.2856  # H122: 3305.7693 -> 3190.2856
$ws.Cells.Item(122, 9).Value = 2147.25  # I122: 2212.7144 -> 2147.25
$ws.Cells.Item(122, 11).Value = 6441.75  # K122: 6638.1432 -> 6441.75
$ws.Cells.Item(122, 13).Value = -3991.75  # M122: -4188.1432 -> -3991.75

$ws.Cells.Item(126, 8).Value = 5385.5713  # H126: 5983.1665 -> 5385.5713
$ws.Cells.Item(126, 9).Value = 5850  # I126: 9900 -> 5850
$ws.Cells.Item(126, 11).Value = 17550  # K126: 29700 -> 17550
$ws.Cells.Item(126, 13).Value = -15080  # M126: -27230 -> -15080

$ws.Cells.Item(132, 8).Value = 8238.9  # H132: 9713.714 -> 8238.9
$ws.Cells.Item(132, 9).Value = 8997.799999999999  # I132: 10332 -> 8997.799999999999
$ws.Cells.Item(132, 10).Value = 7480  # J132: 9250 -> 7480
$ws.Cells.Item(132, 11).Value = 26993.4  # K132: 30996 -> 26993.4
$ws.Cells.Item(132, 12).Value = 22440  # L132: 27750 -> 22440
$ws.Cells.Item(132, 13).Value = -24463.4  # M132: -28466 -> -24463.4
$ws.Cells.Item(132, 14).Value = -27500  # N132: -32810 -> -27500

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 15206.125  # H126: 17178.428 -> 15206.125
$ws.Cells.Item(126, 9).Value = 17021.285  # I126: 19624.834 -> 17021.285
$ws.Cells.Item(126, 11).Value = 51063.855  # K126: 58874.50199999999 -> 51063.855
$ws.Cells.Item(126, 13).Value = -48593.855  # M126: -56404.50199999999 -> -48593.855
